$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.101.76"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.105.40"
$ws.Range("E3").Value = "  -0.44%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "348.74"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("E6").Value = "  -0.31%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5160"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.50%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4442"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "52.65"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.45%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.08938"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.173"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "25.68"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.58%  "
$ws.Range("D13").Value = "2.114.76"
$ws.Range("E13").Value = "  +0.23%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "8.198"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.97%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.728"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "99.01"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.00%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001145"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("E18").Value = "  -0.17%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "20.76"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +7.14%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.06686"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.228"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "30.216.16"
$ws.Range("E23").Value = "  -1.36%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "12.82"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.351"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "2.357.81"
$ws.Range("E26").Value = "  +0.05%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.94"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.80%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.533"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.53%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "162.12"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.93%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "133.50"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.172"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.14%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.1066"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.629"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.82%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.243"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.77%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.974"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "10.39"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.906"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02576"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.39%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.06823"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.2303"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "12.58"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.6818"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.66%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.281"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.77%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "14.30"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.308"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6366"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000364"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.48%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.652"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.221"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.42%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "82.51"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.07237"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
